# Add three time-varying covariate columns (WT, CRCL, DIAL) to the
# NONMEM-style dataset sheet, matching the existing "." missing-value
# convention already used in column D (EVID's neighbour).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New header cells in row 1
$ws.Range("G1").Value = "WT"
$ws.Range("H1").Value = "CRCL"
$ws.Range("I1").Value = "DIAL"

# Placeholder "." values for the single data row (row 2)
$ws.Range("G2").Value = "."
$ws.Range("H2").Value = "."
$ws.Range("I2").Value = "."

# Leave the selection where the author left off editing
$ws.Range("I8").Select()
